# Limit Management Slider Update work
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = "SELECT TM.FIRST_EXECUTION_DATE,TM.LAST_EXECUTION_DATE FROM DC_SCHEDULED_TRAN_MASTER TM WHERE TM.BILL_BENEFICIARY_ID = (SELECT PB.BENEFICIARY_ID FROM DC_BILL_PAYMENT_BENEFICIARY PB WHERE PB.CONSUMER_NUMBER = '{ConsumerNo}' AND PB.CUSTOMER_INFO_ID = (SELECT CI.CUSTOMER_INFO_ID FROM DC_CUSTOMER_INFO CI WHERE CI.CUSTOMER_NAME = '{customer_name}') AND PB.IS_ACTIVE = 1)"

$ws.Range("AA2").Value = $newQuery
$ws.Range("AA3").Value = $newQuery
$ws.Range("AA4").Value = $newQuery
$ws.Range("AA5").Value = $newQuery

$ws.Range("AA1:AA27").ColumnWidth = 203.28515625

$ws.Range("AA10").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 27
$excel.ActiveWindow.ScrollRow = 1
